# Auto-generated edit script applying market-data refresh values
# from the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1475
$ws.Range("I18").Value = 1475
$ws.Range("K18").Value = 1475
$ws.Range("M18").Value = -1191
$ws.Range("H39").Value = 1432.2858
$ws.Range("I39").Value = 73.2
$ws.Range("J39").Value = 4830
$ws.Range("K39").Value = 219.6
$ws.Range("L39").Value = 14490
$ws.Range("M39").Value = 76.39999999999998
$ws.Range("N39").Value = -15082
$ws.Range("H138").Value = 3176.6897
$ws.Range("J138").Value = 3203.0144
$ws.Range("L138").Value = 9609.0432
$ws.Range("N138").Value = -19889.0432
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 479076.53
$ws.Range("I2").Value = 927217.9
$ws.Range("J2").Value = 4573.9414
$ws.Range("K2").Value = 927217.9
$ws.Range("L2").Value = 4573.9414
$ws.Range("M2").Value = -927104.9
$ws.Range("N2").Value = -4799.9414
$ws.Range("H97").Value = 405.14285
$ws.Range("I97").Value = 446.08334
$ws.Range("K97").Value = 446.08334
$ws.Range("M97").Value = 49.91665999999998
$ws.Range("H110").Value = 69059.734
$ws.Range("I110").Value = 78837.16
$ws.Range("K110").Value = 78837.16
$ws.Range("M110").Value = -76792.16
$ws.Range("H116").Value = 479076.53
$ws.Range("I116").Value = 927217.9
$ws.Range("J116").Value = 4573.9414
$ws.Range("K116").Value = 927217.9
$ws.Range("L116").Value = 4573.9414
$ws.Range("M116").Value = -924923.9
$ws.Range("N116").Value = -9161.9414
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 479076.53
$ws.Range("I3").Value = 927217.9
$ws.Range("J3").Value = 4573.9414
$ws.Range("K3").Value = 927217.9
$ws.Range("L3").Value = 4573.9414
$ws.Range("M3").Value = -927103.9
$ws.Range("N3").Value = -4801.9414
$ws.Range("H134").Value = 13159963
$ws.Range("I134").Value = 13159963
$ws.Range("K134").Value = 39479889
$ws.Range("M134").Value = -39477354
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3664.0625
$ws.Range("I86").Value = 3794.9285
$ws.Range("J86").Value = 2748
$ws.Range("K86").Value = 3794.9285
$ws.Range("L86").Value = 2748
$ws.Range("M86").Value = -2671.9285
$ws.Range("N86").Value = -4994
$ws.Range("H89").Value = 3664.0625
$ws.Range("I89").Value = 3794.9285
$ws.Range("J89").Value = 2748
$ws.Range("K89").Value = 18974.6425
$ws.Range("L89").Value = 13740
$ws.Range("M89").Value = -13358.6425
$ws.Range("N89").Value = -24972
$ws.Range("H99").Value = 4831.5
$ws.Range("I99").Value = 5186
$ws.Range("K99").Value = 5186
$ws.Range("M99").Value = -3688
$ws.Range("H105").Value = 5001548
$ws.Range("I105").Value = 5001548
$ws.Range("K105").Value = 5001548
$ws.Range("M105").Value = -4999801
$ws.Range("H107").Value = 670978
$ws.Range("I107").Value = 724873.7
$ws.Range("K107").Value = 724873.7
$ws.Range("M107").Value = -722953.7
$ws.Range("H126").Value = 4831.5
$ws.Range("I126").Value = 5186
$ws.Range("K126").Value = 15558
$ws.Range("M126").Value = -13088
$ws.Range("H131").Value = 98731.25
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = ""
$ws.Range("H132").Value = 111169380
$ws.Range("I132").Value = 142930910
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 428792730
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -428790200
$ws.Range("N132").Value = -17058.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 7.875
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 8.714286
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 52.28571599999999
$ws.Range("M2").Value = 101
$ws.Range("N2").Value = -278.285716
$ws.Range("H6").Value = 66.59999999999999
$ws.Range("I6").Value = 33.25
$ws.Range("K6").Value = 99.75
$ws.Range("M6").Value = 13.25
$ws.Range("H21").Value = 2578.4
$ws.Range("I21").Value = 1995.5
$ws.Range("J21").Value = 2967
$ws.Range("K21").Value = 5986.5
$ws.Range("L21").Value = 8901
$ws.Range("M21").Value = -5813.5
$ws.Range("N21").Value = -9247
$ws.Range("H128").Value = 194849.5
$ws.Range("I128").Value = 194849.5
$ws.Range("K128").Value = 584548.5
$ws.Range("M128").Value = -579568.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25629.25
$ws.Range("J15").Value = 26290.572
$ws.Range("L15").Value = 26290.572
$ws.Range("N15").Value = -26866.572
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H80").Value = 2988.5334
$ws.Range("I80").Value = 2204.2856
$ws.Range("J80").Value = 3674.75
$ws.Range("K80").Value = 2204.2856
$ws.Range("L80").Value = 3674.75
$ws.Range("M80").Value = -1206.2856
$ws.Range("N80").Value = -5670.75
$ws.Range("H81").Value = 25629.25
$ws.Range("J81").Value = 26290.572
$ws.Range("L81").Value = 26290.572
$ws.Range("N81").Value = -28286.572
$ws.Range("H82").Value = 59994.5
$ws.Range("J82").Value = 59999
$ws.Range("L82").Value = 59999
$ws.Range("N82").Value = -60765
$ws.Range("H83").Value = 2988.5334
$ws.Range("I83").Value = 2204.2856
$ws.Range("J83").Value = 3674.75
$ws.Range("K83").Value = 11021.428
$ws.Range("L83").Value = 18373.75
$ws.Range("M83").Value = -6029.428
$ws.Range("N83").Value = -28357.75
$ws.Range("H84").Value = 25629.25
$ws.Range("J84").Value = 26290.572
$ws.Range("L84").Value = 78871.716
$ws.Range("N84").Value = -88855.716
$ws.Range("H85").Value = 59994.5
$ws.Range("J85").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("N85").Value = -62651
$ws.Range("H113").Value = 132309.25
$ws.Range("I113").Value = 205894.8
$ws.Range("K113").Value = 205894.8
$ws.Range("M113").Value = -203724.8
$ws.Range("H132").Value = 2453197
$ws.Range("I132").Value = 2453197
$ws.Range("K132").Value = 7359591
$ws.Range("M132").Value = -7357061
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 29995
$ws.Range("J105").Value = 29995
$ws.Range("L105").Value = 29995
$ws.Range("N105").Value = -36983
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5370.8237
$ws.Range("I122").Value = 5370.8237
$ws.Range("K122").Value = 16112.4711
$ws.Range("M122").Value = -13662.4711
$ws.Range("H126").Value = 2850.5715
$ws.Range("I126").Value = 2850.5715
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8551.7145
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -6081.7145
$ws.Range("H136").Value = 41670064
$ws.Range("I136").Value = 62502124
$ws.Range("K136").Value = 187506372
$ws.Range("M136").Value = -187503822
